# Refresh the crypto price/volume snapshot (GitHub Actions data pull).
# Source cells are plain text (t="inlineStr"/shared-string), never numbers,
# so numeric-looking Price values are written with a leading apostrophe
# (Excel's quote-prefix) to force text storage instead of silently
# converting them to floating point numbers and losing formatting
# like trailing zeros (e.g. "68.00") or multi-dot groupings (e.g. "41.528.77").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "41.528.77"
$ws.Range("E2").Value = "  +0.59%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value = "2.476.92"
$ws.Range("E3").Value = "  +0.59%  "

# Row 4 (TetherUSD)
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.23%  "

# Row 5 (BNB)
$ws.Range("D5").Value = "'313.23"
$ws.Range("E5").Value = "  +0.29%  "

# Row 6 (Solana)
$ws.Range("D6").Value = "'92.56"
$ws.Range("E6").Value = "  -1.69%  "

# Row 7 (XRP)
$ws.Range("E7").Value = "  -1.18%  "

# Row 8 (USDC)
$ws.Range("D8").Value = "'0.999"

# Row 9 (Cardano)
$ws.Range("D9").Value = "'0.505"
$ws.Range("E9").Value = "  +1.50%  "

# Row 10 (Avalanche)
$ws.Range("D10").Value = "'32.68"
$ws.Range("E10").Value = "  -2.20%  "

# Row 11 (Dogecoin)
$ws.Range("D11").Value = "'0.0786"
$ws.Range("E11").Value = "  +0.89%  "

# Row 12 (TRON)
$ws.Range("D12").Value = "'0.111"
$ws.Range("E12").Value = "  +2.20%  "

# Row 13 (WrappedliquidstakedEther2.0)
$ws.Range("D13").Value = "2.858.71"
$ws.Range("E13").Value = "  +0.57%  "

# Row 14 (Chainlink)
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'6.85"
$ws.Range("E14").Value = "  -2.16%  "

# Row 15 (Polkadot)
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'16.14"
$ws.Range("E15").Value = "  +8.47%  "

# Row 16 (WrappedEther)
$ws.Range("D16").Value = "2.473.70"
$ws.Range("E16").Value = "  +0.49%  "

# Row 17 (Polygon)
$ws.Range("D17").Value = "'0.767"
$ws.Range("E17").Value = "  -2.09%  "

# Row 18 (WrappedBTC)
$ws.Range("D18").Value = "41.516.15"
$ws.Range("E18").Value = "  +0.71%  "

# Row 19 (Uniswap)
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0950"
$ws.Range("E19").Value = "  +3.08%  "

# Row 20 (ShibaInu)
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'6.46"
$ws.Range("E20").Value = "  +2.44%  "

# Row 21 (Litecoin)
$ws.Range("D21").Value = "'72.23"
$ws.Range("E21").Value = "  +5.49%  "

# Row 22 (InternetComputer(DFINITY))
$ws.Range("D22").Value = "'11.17"
$ws.Range("E22").Value = "  -0.84%  "

# Row 23 (BitcoinCash)
$ws.Range("D23").Value = "'236.19"
$ws.Range("E23").Value = "  -0.18%  "

# Row 24 (PancakeSwap)
$ws.Range("E24").Value = "  -1.82%  "

# Row 25 (Dai)
$ws.Range("E25").Value = "  -0.14%  "

# Row 26 (ImmutableX)
$ws.Range("D26").Value = "'1.90"
$ws.Range("E26").Value = "  -0.41%  "

# Row 27 (EthereumClassic)
$ws.Range("D27").Value = "'24.76"
$ws.Range("E27").Value = "  +2.88%  "

# Row 28 (Toncoin)
$ws.Range("E28").Value = "  +0.12%  "

# Row 29 (Cosmos)
$ws.Range("D29").Value = "'9.63"
$ws.Range("E29").Value = "  +0.00%  "

# Row 30 (InjectiveProtocol)
$ws.Range("D30").Value = "'35.76"
$ws.Range("E30").Value = "  -2.64%  "

# Row 31 (Monero)
$ws.Range("D31").Value = "'158.07"
$ws.Range("E31").Value = "  +3.83%  "

# Row 32 (Filecoin)
$ws.Range("E32").Value = "  -1.12%  "

# Row 33 (WEMIXToken)
$ws.Range("D33").Value = "'2.57"
$ws.Range("E33").Value = "  +0.38%  "

# Row 34 (Hedera)
$ws.Range("D34").Value = "'0.0757"
$ws.Range("E34").Value = "  +1.78%  "

# Row 35 (Celestia)
$ws.Range("D35").Value = "'17.36"
$ws.Range("E35").Value = "  +2.14%  "

# Row 36 (ApeXProtocol)
$ws.Range("D36").Value = "'2.41"
$ws.Range("E36").Value = "  -9.21%  "

# Row 37 (Kaspa)
$ws.Range("D37").Value = "'0.107"
$ws.Range("E37").Value = "  +4.42%  "

# Row 38 (LidoDAOToken)
$ws.Range("E38").Value = "  -4.82%  "

# Row 39 (ARBITRUM)
$ws.Range("E39").Value = "  -3.62%  "

# Row 40 (Stellar)
$ws.Range("E40").Value = "  +0.07%  "

# Row 41 (RenderToken)
$ws.Range("D41").Value = "'4.10"
$ws.Range("E41").Value = "  -3.62%  "

# Row 42 (FirstDigitalUSD)
$ws.Range("E42").Value = "  -0.29%  "

# Row 43 (Maker)
$ws.Range("D43").Value = "1.984.82"
$ws.Range("E43").Value = "  -0.02%  "

# Row 44 (EnergySwap)
$ws.Range("D44").Value = "'19.24"
$ws.Range("E44").Value = "  -3.36%  "

# Row 45 (VeChain)
$ws.Range("E45").Value = "  -0.31%  "

# Row 46 (NEARProtocol)
$ws.Range("E46").Value = "  -3.14%  "

# Row 47 (FraxShare)
$ws.Range("D47").Value = "'8.97"
$ws.Range("E47").Value = "  +1.92%  "

# Row 48 (RocketPoolETH)
$ws.Range("D48").Value = "2.717.45"
$ws.Range("E48").Value = "  +0.47%  "

# Row 49 (Aave)
$ws.Range("D49").Value = "'97.32"
$ws.Range("E49").Value = "  +0.28%  "

# Row 50 (ordi)
$ws.Range("D50").Value = "'68.00"
$ws.Range("E50").Value = "  -1.48%  "

# Row 51 (BitcoinSV)
$ws.Range("D51").Value = "'72.29"
$ws.Range("E51").Value = "  -2.80%  "
